# Second table "成分-UTS（铸铝）" (Composition-UTS, cast aluminum) is added to the
# right of the existing "成分-UTS（所有可用数据）" table, reusing the same
# model/R2/MAPE layout. RF/KNN/SVR rows are not computed yet, so they are
# copied over from the first table as placeholders (per commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# --- New merged title cell D2:F2 -------------------------------------------------
$ws.Range("D2").Value = "成分-UTS（铸铝）"
$ws.Range("D2:F2").Merge()
$ws.Range("D2:F2").HorizontalAlignment = $xlCenter

# --- New header row D3:F3 (模型 / R2 / MAPE) -------------------------------------
$ws.Cells.Item(3, 4).Value = "模型"
$ws.Cells.Item(3, 5).Value = "R2"
$ws.Cells.Item(3, 6).Value = "MAPE"
$ws.Range("D3:F3").HorizontalAlignment = $xlCenter

# --- New data rows D4:F9 ---------------------------------------------------------
# model name, R2, MAPE
$data = @(
    @("LIN",  0.5158, 0.1659),
    @("POLY", 0.6982, 0.1373),
    @("DT",   0.1969, 0.2029),
    @("RF",   0.7433, 0.1553),
    @("KNN",  0.7155, 0.1645),
    @("SVR",  0.7226, 0.1558)
)

$row = 4
foreach ($item in $data) {
    $ws.Cells.Item($row, 4).Value = $item[0]
    $ws.Cells.Item($row, 5).Value = $item[1]
    $ws.Cells.Item($row, 6).Value = $item[2]
    $row++
}

$ws.Range("D4:D9").HorizontalAlignment = $xlCenter
$ws.Range("E4:F9").HorizontalAlignment = $xlCenter
$ws.Range("E4:F9").NumberFormat = "0.00%"

# --- Restore selection like in the saved workbook --------------------------------
$ws.Range("F6").Select()
